$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 392 (shifts old rows 392.. down to 394..)
$ws.Range("A392:T393").EntireRow.Insert()

# New row 392
$ws.Range("A392").Value = 7
$ws.Range("B392").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C392").Value = "Ñuble"
$ws.Range("D392").Value = 44505
$ws.Range("E392").Value = 16
$ws.Range("F392").Value = "Fruta"
$ws.Range("G392").Value = 100102
$ws.Range("H392").Value = "Cítricos"
$ws.Range("I392").Value = 100102003
$ws.Range("J392").Value = "Limón"
$ws.Range("K392").Value = "Sin especificar"
$ws.Range("L392").Value = "1a amarillo"
$ws.Range("M392").Value = 160
$ws.Range("N392").Value = 5800
$ws.Range("O392").Value = 6000
$ws.Range("P392").Value = 5900
$ws.Range("Q392").Value = '$/malla 16 kilos'
$ws.Range("R392").Value = "Región de O'Higgins"
$ws.Range("S392").Value = 369
$ws.Range("T392").Value = 16

# New row 393
$ws.Range("A393").Value = 7
$ws.Range("B393").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C393").Value = "Ñuble"
$ws.Range("D393").Value = 44505
$ws.Range("E393").Value = 16
$ws.Range("F393").Value = "Fruta"
$ws.Range("G393").Value = 100102
$ws.Range("H393").Value = "Cítricos"
$ws.Range("I393").Value = 100102003
$ws.Range("J393").Value = "Limón"
$ws.Range("K393").Value = "Sin especificar"
$ws.Range("L393").Value = "2a amarillo"
$ws.Range("M393").Value = 160
$ws.Range("N393").Value = 4500
$ws.Range("O393").Value = 5000
$ws.Range("P393").Value = 4750
$ws.Range("Q393").Value = '$/malla 16 kilos'
$ws.Range("R393").Value = "Región de O'Higgins"
$ws.Range("S393").Value = 297
$ws.Range("T393").Value = 16

# Ensure date cells keep the date number format (style index 2) as used elsewhere in column D
$ws.Range("D392").NumberFormat = $ws.Range("D391").NumberFormat
$ws.Range("D393").NumberFormat = $ws.Range("D391").NumberFormat
